# Generate Report for Handback
#
# The earlier handback run (GUID1, d93a9d40-... -> 438699ce-...) is re-synced
# to its newer timestamps/hashes, and a second handback file (GUID2,
# dc90e2a8-...) is appended as a new row to all three sheets (Overview,
# zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$GUID1 = "438699ce-5cb5-4e86-822d-2b503fedfca4"
$GUID2 = "dc90e2a8-89e2-4731-b5ab-66a7f86b0c93"
$HASH1 = "7bdecfbf4e8cacf10ad7deeb9e1485049a6645be"
$HASH2 = "6faa6f6cf1fce11535fd9e4732dfea1382638e2b"

# Leading apostrophe forces literal text so values that look like booleans
# ("True"/"False") or blanks don't silently become t="b"/omitted cells.
function T([string]$s) { return "'" + $s }

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item("Overview")

# Hyperlinks must be rebuilt as a batch: deleting via a single-cell Range's
# Hyperlinks collection clears the whole sheet's collection in this engine.
$wsO.Hyperlinks.Delete()

# Update existing row 2 (re-run of the same source file) to the later timestamp.
$wsO.Range("A2").Value2 = (T "$GUID1.md")
$wsO.Range("B2").Value2 = (T "e2e\$GUID1.md")
$wsO.Range("C2").Value2 = (T ".md")
$wsO.Range("E2").Value2 = (T "Handed back: in sync with en-US")
$wsO.Range("F2").Value2 = (T "Handed back: in sync with en-US")
$wsO.Range("G2").Value2 = (T "2016-08-27 16:59:28")

# Append row 3 for the new handback file.
$loO.ListRows.Add() | Out-Null
$wsO.Range("A3").Value2 = (T "$GUID2.md")
$wsO.Range("B3").Value2 = (T "e2e\$GUID2.md")
$wsO.Range("C3").Value2 = (T ".md")
$wsO.Range("E3").Value2 = (T "Handed back: in sync with en-US")
$wsO.Range("F3").Value2 = (T "Handed back: in sync with en-US")
$wsO.Range("G3").Value2 = (T "2016-08-27 16:59:28")

$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID1.md", [Type]::Missing, [Type]::Missing, "e2e\$GUID1.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID2.md", [Type]::Missing, [Type]::Missing, "e2e\$GUID2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item("zh-cn")

$wsZ.Hyperlinks.Delete()

$wsZ.Range("A2").Value2 = (T "$GUID1.md")
$wsZ.Range("B2").Value2 = (T ".md")
$wsZ.Range("C2").Value2 = (T "Handed back: in sync with en-US")
$wsZ.Range("D2").Value2 = (T "e2e")
$wsZ.Range("E2").Value2 = (T "ht")
$wsZ.Range("F2").Value2 = (T "False")
$wsZ.Range("G2").Value2 = (T "$GUID1.$HASH1.zh-cn.xlf")
$wsZ.Range("H2").Value2 = (T "2016-08-27 16:59:23")
$wsZ.Range("I2").Value2 = (T "$GUID1.md")
$wsZ.Range("J2").Value2 = (T "$GUID1.$HASH1.zh-cn.xlf")
$wsZ.Range("K2").Value2 = (T "2016-08-27 16:59:40")
$wsZ.Range("L2").Value2 = (T "")
$wsZ.Range("M2").Value2 = (T "True")
$wsZ.Range("N2").Value2 = (T "")
$wsZ.Range("O2").Value2 = (T "False")
$wsZ.Range("P2").Value2 = (T "")

$loZ.ListRows.Add() | Out-Null
$wsZ.Range("A3").Value2 = (T "$GUID2.md")
$wsZ.Range("B3").Value2 = (T ".md")
$wsZ.Range("C3").Value2 = (T "Handed back: in sync with en-US")
$wsZ.Range("D3").Value2 = (T "e2e")
$wsZ.Range("E3").Value2 = (T "ht")
$wsZ.Range("F3").Value2 = (T "True")
$wsZ.Range("G3").Value2 = (T "$GUID2.$HASH2.zh-cn.xlf")
$wsZ.Range("H3").Value2 = (T "2016-08-27 16:59:23")
$wsZ.Range("I3").Value2 = (T "$GUID2.md")
$wsZ.Range("J3").Value2 = (T "$GUID2.$HASH2.zh-cn.xlf")
$wsZ.Range("K3").Value2 = (T "2016-08-27 16:59:40")
$wsZ.Range("L3").Value2 = (T "")
$wsZ.Range("M3").Value2 = (T "True")
$wsZ.Range("N3").Value2 = (T "")
$wsZ.Range("O3").Value2 = (T "False")
$wsZ.Range("P3").Value2 = (T "")

$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID1.md", [Type]::Missing, [Type]::Missing, "$GUID1.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/dc1179276f9a909992490ffd2c675927138f2648/e2e/$GUID1.md", [Type]::Missing, [Type]::Missing, "$GUID1.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID2.md", [Type]::Missing, [Type]::Missing, "$GUID2.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/dc1179276f9a909992490ffd2c675927138f2648/e2e/$GUID2.md", [Type]::Missing, [Type]::Missing, "$GUID2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item("de-de")

$wsD.Hyperlinks.Delete()

$wsD.Range("A2").Value2 = (T "$GUID1.md")
$wsD.Range("B2").Value2 = (T ".md")
$wsD.Range("C2").Value2 = (T "Handed back: in sync with en-US")
$wsD.Range("D2").Value2 = (T "e2e")
$wsD.Range("E2").Value2 = (T "ht")
$wsD.Range("F2").Value2 = (T "False")
$wsD.Range("G2").Value2 = (T "$GUID1.$HASH1.de-de.xlf")
$wsD.Range("H2").Value2 = (T "2016-08-27 16:59:28")
$wsD.Range("I2").Value2 = (T "$GUID1.md")
$wsD.Range("J2").Value2 = (T "$GUID1.$HASH1.de-de.xlf")
$wsD.Range("K2").Value2 = (T "2016-08-27 16:59:47")
$wsD.Range("L2").Value2 = (T "")
$wsD.Range("M2").Value2 = (T "True")
$wsD.Range("N2").Value2 = (T "")
$wsD.Range("O2").Value2 = (T "False")
$wsD.Range("P2").Value2 = (T "")

$loD.ListRows.Add() | Out-Null
$wsD.Range("A3").Value2 = (T "$GUID2.md")
$wsD.Range("B3").Value2 = (T ".md")
$wsD.Range("C3").Value2 = (T "Handed back: in sync with en-US")
$wsD.Range("D3").Value2 = (T "e2e")
$wsD.Range("E3").Value2 = (T "ht")
$wsD.Range("F3").Value2 = (T "True")
$wsD.Range("G3").Value2 = (T "$GUID2.$HASH2.de-de.xlf")
$wsD.Range("H3").Value2 = (T "2016-08-27 16:59:28")
$wsD.Range("I3").Value2 = (T "$GUID2.md")
$wsD.Range("J3").Value2 = (T "$GUID2.$HASH2.de-de.xlf")
$wsD.Range("K3").Value2 = (T "2016-08-27 16:59:47")
$wsD.Range("L3").Value2 = (T "")
$wsD.Range("M3").Value2 = (T "True")
$wsD.Range("N3").Value2 = (T "")
$wsD.Range("O3").Value2 = (T "False")
$wsD.Range("P3").Value2 = (T "")

$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID1.md", [Type]::Missing, [Type]::Missing, "$GUID1.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6f87c9e2552646c596a2828c89c708468205fd86/e2e/$GUID1.md", [Type]::Missing, [Type]::Missing, "$GUID1.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9a99c37fadaec91ddbc5ec55a14291a337b5bde/e2e/$GUID2.md", [Type]::Missing, [Type]::Missing, "$GUID2.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6f87c9e2552646c596a2828c89c708468205fd86/e2e/$GUID2.md", [Type]::Missing, [Type]::Missing, "$GUID2.md") | Out-Null

Write-Host "Done applying handback report update."
